$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "日期：2021/12/22"
$ws.Range("B2").Value = "202202"
$ws.Range("C2").Value = 17802
$ws.Range("D2").Value = 2496
$ws.Range("E2").Value = 605268
$ws.Range("F2").Value = 17638

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "日期：2021/12/22"
$ws.Range("B2").Value = 0.08

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "110年12月22日"
$ws.Range("B2").Value = 17.7
$ws.Range("C2").Value = 19.86

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "110年12月22日"
$ws.Range("B2").Value = 17835.4

# --- Sheet 5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "2021/12/22"
$ws.Range("B2").Value = 46323
$ws.Range("C2").Value = 53049
$ws.Range("D2").Value = -1116
$ws.Range("E2").Value = -2269
$ws.Range("F2").Value = 21197
$ws.Range("G2").Value = 45391
$ws.Range("H2").Value = -909
$ws.Range("I2").Value = -2128
$ws.Range("J2").Value = -24194
$ws.Range("K2").Value = 1219
$ws.Range("L2").Value = -207
$ws.Range("M2").Value = -141
$ws.Range("N2").Value = -66

# --- Sheet1: materialize trailing empty row (dimension -> A1:F8) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(8, 1).Borders.LineStyle = -4142
